# Adapt column header formatting to respective input file names (#7)
# - rename "<Column>_old" -> "<Column>_FV2410" and "<Column>_new" -> "<Column>_FV2504"
# - turn the sheet's range into a real Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header cells (A1:J1 = "..._old" columns, L1:U1 = "..._new"
#    columns; K1 stays "diff").
# ---------------------------------------------------------------------------
$oldHeaders = @(
  "Segmentname_old",
  "Segmentgruppe_old",
  "Segment_old",
  "Datenelement_old",
  "Segment ID_old",
  "Code_old",
  "Qualifier_old",
  "Beschreibung_old",
  "Bedingungsausdruck_old",
  "Bedingung_old"
)
$newHeaders = @(
  "Segmentname_new",
  "Segmentgruppe_new",
  "Segment_new",
  "Datenelement_new",
  "Segment ID_new",
  "Code_new",
  "Qualifier_new",
  "Beschreibung_new",
  "Bedingungsausdruck_new",
  "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace '_old$', '_FV2410')
}
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace '_new$', '_FV2504')
}

# ---------------------------------------------------------------------------
# 2) Convert A1:U66 into a native Excel Table ("Table1") without letting the
#    table-creation step capture the header row's pre-existing formatting as
#    a new dxf (the original workbook keeps dxfs count="0"). We stash the
#    header formatting on an unused scratch row, strip the header formatting
#    so table creation has nothing to snapshot, build the table, then copy
#    the formatting back onto the header row and scrub the scratch row.
# ---------------------------------------------------------------------------
$headerRange  = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$headerRange.Copy()
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U66"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$scratchRange.ClearFormats()
$scratchRange.ClearContents()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
